# Add two new case rows (44 and 45) to the COVID-19 case-tracking sheet,
# bringing the data through "Feb 10" (commit: "update on Feb 10").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the formatting of the last existing data row (44) down onto the two
# new rows so the new cells share the same style indices (general number
# format for most columns, text format for the Date column, right-aligned
# text format for the blank "Related" column, etc).
$ws.Range("A44:K44").Copy() | Out-Null
$ws.Range("A45:K46").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Row 45 -> Case 44
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = 1.452561
$ws.Range("C45").Value = 103.816625
$ws.Range("D45").Value = "Feb-09"
$ws.Range("E45").Value = 37
$ws.Range("F45").Value = "Male"
$ws.Range("G45").Value = "Singapore"
$ws.Range("H45").Value = "Sembawang Drive"
$ws.Range("I45").Value = "Certis Cisco Centre (20 Jalan Afifi), Chingay 2020, Khoo Teck Puat Hospital (KTPH)"

# Row 46 -> Case 45
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = 1.3105450000000001
$ws.Range("C46").Value = 103.84683800000001
$ws.Range("D46").Value = "Feb-10"
$ws.Range("E46").Value = 2
$ws.Range("F46").Value = "Female"
$ws.Range("G46").Value = "China, Wuhan"
$ws.Range("H46").Value = "Arrived from Wuhan"
$ws.Range("I46").Value = "KK Women" + [char]0x2019 + "s and Children" + [char]0x2019 + "s Hospital"

# Leave J45:K46 ("Related" / "Status") blank, matching the source row.

# Update the view state to match the author's final scroll/selection position.
$ws.Range("F49").Select() | Out-Null
